$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their text formatting so values
# like "1.001" or "74.72" are not reinterpreted as numbers/dates.

$priceCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.505.85"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "1.666.04"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "233.79"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.4620"
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "0.06127"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "1.666.11"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "0.06947"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "4.332"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "74.72"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "0.5623"
$ws.Range("E15").Value = "  -6.44%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "25.509.00"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("D19").Value = "0.000006676"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "1.879.13"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").Value = "4.411"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "8.691"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "5.186"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").Value = "136.40"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").Value = "14.83"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "1.369"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "104.13"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "1.694"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").Value = "3.933"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "0.07722"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "0.04264"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "2.626"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "0.9408"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").Value = "0.5956"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "0.9205"
$ws.Range("E37").Value = "  +12.49%  "
$ws.Range("D38").Value = "2.476"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").Value = "1.000"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "101.70"
$ws.Range("E40").Value = "  +3.94%  "
$ws.Range("D41").Value = "0.01460"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").Value = "1.809"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "0.3691"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "4.914"
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("D45").Value = "0.05291"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "6.104"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "7.337"
$ws.Range("E51").Value = "  +1.69%  "
